$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.13533147010429
$ws.Range("C2").Value = 0.2252243830160694
$ws.Range("D2").Value = 0.07834349447604438
$ws.Range("E2").Value = 0.09406137471690301
$ws.Range("G2").Value = 0.002455704863333368
$ws.Range("I2").Value = 0.836628255675425
$ws.Range("L2").Value = 0.2081373669001252
$ws.Range("M2").Value = 0.2467649940627155
$ws.Range("O2").Value = 3.497469496676302
# Row 3
$ws.Range("B3").Value = 1.030758460242964
$ws.Range("C3").Value = 0.2047744254066401
$ws.Range("D3").Value = 0.0710865866640944
$ws.Range("E3").Value = 0.09485364372930793
$ws.Range("G3").Value = 0.002459056076272053
$ws.Range("I3").Value = 0.8469054332425046
$ws.Range("L3").Value = 0.2055324021110678
$ws.Range("M3").Value = 0.2305992805664658
$ws.Range("O3").Value = 3.509885920984146
# Row 4
$ws.Range("B4").Value = 0.9666470508117868
$ws.Range("C4").Value = 0.1921354109002209
$ws.Range("D4").Value = 0.06666704055751893
$ws.Range("E4").Value = 0.09536812119174365
$ws.Range("G4").Value = 0.002461223740790418
$ws.Range("I4").Value = 0.8537931240883552
$ws.Range("L4").Value = 0.2040389907457012
$ws.Range("M4").Value = 0.2207373964448536
$ws.Range("O4").Value = 3.520103826522131
# Row 5
$ws.Range("B5").Value = 0.9405468446136069
$ws.Range("C5").Value = 0.1869643173884867
$ws.Range("D5").Value = 0.06487510809674291
$ws.Range("E5").Value = 0.09558483383782779
$ws.Range("G5").Value = 0.002462134829020659
$ws.Range("I5").Value = 0.856744895178629
$ws.Range("L5").Value = 0.2034571138309715
$ws.Range("M5").Value = 0.2167348745354545
$ws.Range("O5").Value = 3.524918802188751
# Row 6
$ws.Range("B6").Value = 0.9362145187471924
$ws.Range("C6").Value = 0.1861044232474285
$ws.Range("D6").Value = 0.06457810561474275
$ws.Range("E6").Value = 0.09562124552398643
$ws.Range("G6").Value = 0.002462287792960371
$ws.Range("I6").Value = 0.8572437843879648
$ws.Range("L6").Value = 0.2033621072553089
$ws.Range("M6").Value = 0.2160712478005351
$ws.Range("O6").Value = 3.525757613852733
# Row 7
$ws.Range("B7").Value = 0.9662949483725924
$ws.Range("C7").Value = 0.19206575484327
$ws.Range("D7").Value = 0.06664283726432529
$ws.Range("E7").Value = 0.09537101525591796
$ws.Range("G7").Value = 0.00246123591573722
$ws.Range("I7").Value = 0.8538323459861239
$ws.Range("L7").Value = 0.204031035196337
$ws.Range("M7").Value = 0.2206833508664587
$ws.Range("O7").Value = 3.520166128347768
# Row 8
$ws.Range("B8").Value = 1.099255376459837
$ws.Range("C8").Value = 0.2181905020754016
$ws.Range("D8").Value = 0.07583376757256133
$ws.Range("E8").Value = 0.09432874102463629
$ws.Range("G8").Value = 0.002456837581079772
$ws.Range("I8").Value = 0.8400518384749311
$ws.Range("L8").Value = 0.2072171837096164
$ws.Range("M8").Value = 0.2411779017418354
$ws.Range("O8").Value = 3.501211512862454
# Row 9
$ws.Range("B9").Value = 1.360712704484627
$ws.Range("C9").Value = 0.2687591891919112
$ws.Range("D9").Value = 0.09414736960326309
$ws.Range("E9").Value = 0.09250658835838665
$ws.Range("G9").Value = 0.002449081418067186
$ws.Range("I9").Value = 0.8176203806963365
$ws.Range("L9").Value = 0.2143057612326373
$ws.Range("M9").Value = 0.2818683441421967
$ws.Range("O9").Value = 3.484684402030268
# Row 10
$ws.Range("B10").Value = 1.553204178449221
$ws.Range("C10").Value = 0.3055040490165197
$ws.Range("D10").Value = 0.1077845431347413
$ws.Range("E10").Value = 0.09130219562391728
$ws.Range("G10").Value = 0.002443907246386589
$ws.Range("I10").Value = 0.8039534648840601
$ws.Range("L10").Value = 0.2200258480551156
$ws.Range("M10").Value = 0.3120633415585843
$ws.Range("O10").Value = 3.485213883558146
# Row 11
$ws.Range("B11").Value = 1.640852806360783
$ws.Range("C11").Value = 0.3221309701734185
$ws.Range("D11").Value = 0.1140292108001688
$ws.Range("E11").Value = 0.09078327896951066
$ws.Range("G11").Value = 0.002441666062314168
$ws.Range("I11").Value = 0.7983499801060319
$ws.Range("L11").Value = 0.2227392605101102
$ws.Range("M11").Value = 0.3258639847361238
$ws.Range("O11").Value = 3.488225547504612
# Row 12
$ws.Range("B12").Value = 1.674054003539993
$ws.Range("C12").Value = 0.3284142919544593
$ws.Range("D12").Value = 0.1163998726040205
$ws.Range("E12").Value = 0.09059093118226702
$ws.Range("G12").Value = 0.002440833483941417
$ws.Range("I12").Value = 0.7963165854913399
$ws.Range("L12").Value = 0.2237827483229324
$ws.Range("M12").Value = 0.3310990989131355
$ws.Range("O12").Value = 3.489765884397599
# Row 13
$ws.Range("B13").Value = 1.666903084166449
$ws.Range("C13").Value = 0.3270616446127974
$ws.Range("D13").Value = 0.1158890432930519
$ws.Range("E13").Value = 0.09063217211030716
$ws.Range("G13").Value = 0.002441012079428026
$ws.Range("I13").Value = 0.7967505718495218
$ws.Range("L13").Value = 0.2235573044225418
$ws.Range("M13").Value = 0.3299712214903607
$ws.Range("O13").Value = 3.4894163372484
# Row 14
$ws.Range("B14").Value = 1.643584085151019
$ws.Range("C14").Value = 0.3226481631940317
$ws.Range("D14").Value = 0.1142241272715978
$ws.Range("E14").Value = 0.09076737119683775
$ws.Range("G14").Value = 0.002441597243198158
$ws.Range("I14").Value = 0.7981809157651725
$ws.Range("L14").Value = 0.2228247888385084
$ws.Range("M14").Value = 0.3262944993226782
$ws.Range("O14").Value = 3.488344249368168
# Row 15
$ws.Range("B15").Value = 1.629301859024224
$ws.Range("C15").Value = 0.3199430881363696
$ws.Range("D15").Value = 0.1132050930699364
$ws.Range("E15").Value = 0.09085072528450144
$ws.Range("G15").Value = 0.002441957768601362
$ws.Range("I15").Value = 0.7990685793725802
$ws.Range("L15").Value = 0.2223781817454267
$ws.Range("M15").Value = 0.3240435811579303
$ws.Range("O15").Value = 3.487739684370808
# Row 16
$ws.Range("B16").Value = 1.547477689330265
$ws.Range("C16").Value = 0.3044156420764637
$ws.Range("D16").Value = 0.1073772687552577
$ws.Range("E16").Value = 0.09133669002324485
$ws.Range("G16").Value = 0.002444055971587031
$ws.Range("I16").Value = 0.8043320381556285
$ws.Range("L16").Value = 0.2198507575879631
$ws.Range("M16").Value = 0.3111627230756469
$ws.Range("O16").Value = 3.48507294357313
# Row 17
$ws.Range("B17").Value = 1.497301533524308
$ws.Range("C17").Value = 0.2948672422465108
$ws.Range("D17").Value = 0.1038126286549925
$ws.Range("E17").Value = 0.09164222522631804
$ws.Range("G17").Value = 0.002445371927835449
$ws.Range("I17").Value = 0.807718372412964
$ws.Range("L17").Value = 0.2183287542294607
$ws.Range("M17").Value = 0.3032771799713103
$ws.Range("O17").Value = 3.484147683102663
# Row 18
$ws.Range("B18").Value = 1.468449458956002
$ws.Range("C18").Value = 0.2893669373313799
$ws.Range("D18").Value = 0.1017662033681148
$ws.Range("E18").Value = 0.09182068820664613
$ws.Range("G18").Value = 0.002446139431547385
$ws.Range("I18").Value = 0.8097238465390433
$ws.Range("L18").Value = 0.2174638175255836
$ws.Range("M18").Value = 0.2987477344177307
$ws.Range("O18").Value = 3.483876212079423
# Row 19
$ws.Range("B19").Value = 1.458682046394642
$ws.Range("C19").Value = 0.2875032063227252
$ws.Range("D19").Value = 0.1010739809432408
$ws.Range("E19").Value = 0.09188158143822811
$ws.Range("G19").Value = 0.00244640111810589
$ws.Range("I19").Value = 0.8104127747392624
$ws.Range("L19").Value = 0.2171727655731814
$ws.Range("M19").Value = 0.2972151977482298
$ws.Range("O19").Value = 3.483829030490625
# Row 20
$ws.Range("B20").Value = 1.502642061242454
$ws.Range("C20").Value = 0.2958845480822561
$ws.Range("D20").Value = 0.1041916907561671
$ws.Range("E20").Value = 0.0916094183049283
$ws.Range("G20").Value = 0.002445230745688661
$ws.Range("I20").Value = 0.8073519130125888
$ws.Range("L20").Value = 0.21848968972634
$ws.Range("M20").Value = 0.304115978607733
$ws.Range("O20").Value = 3.484219185203216
# Row 21
$ws.Range("B21").Value = 1.650433164897777
$ws.Range("C21").Value = 0.3239448622825591
$ws.Range("D21").Value = 0.1147129919001202
$ws.Range("E21").Value = 0.090727547276507
$ws.Range("G21").Value = 0.002441424929775983
$ws.Range("I21").Value = 0.7977583843571026
$ws.Range("L21").Value = 0.22303951303536
$ws.Range("M21").Value = 0.3273741961343362
$ws.Range("O21").Value = 3.488648283267736
# Row 22
$ws.Range("B22").Value = 1.747084125850108
$ws.Range("C22").Value = 0.3422084260816689
$ws.Range("D22").Value = 0.1216239210525316
$ws.Range("E22").Value = 0.09017540532555612
$ws.Range("G22").Value = 0.002439031471250544
$ws.Range("I22").Value = 0.7920045709513701
$ws.Range("L22").Value = 0.2261062005668464
$ws.Range("M22").Value = 0.3426277925129071
$ws.Range("O22").Value = 3.493874462476185
# Row 23
$ws.Range("B23").Value = 1.695494576386238
$ws.Range("C23").Value = 0.3324677915951213
$ws.Range("D23").Value = 0.1179322447501505
$ws.Range("E23").Value = 0.09046788205117062
$ws.Range("G23").Value = 0.002440300342501309
$ws.Range("I23").Value = 0.7950281728445674
$ws.Range("L23").Value = 0.2244609414258747
$ws.Range("M23").Value = 0.3344818784625048
$ws.Range("O23").Value = 3.490871341188665
# Row 24
$ws.Range("B24").Value = 1.50022762454978
$ws.Range("C24").Value = 0.2954246577509423
$ws.Range("D24").Value = 0.104020307646806
$ws.Range("E24").Value = 0.09162424157059057
$ws.Range("G24").Value = 0.00244529454003846
$ws.Range("I24").Value = 0.807517406708584
$ws.Range("L24").Value = 0.2184168993783828
$ws.Range("M24").Value = 0.3037367450885071
$ws.Range("O24").Value = 3.484186047834669
# Row 25
$ws.Range("B25").Value = 1.289908602845685
$ws.Range("C25").Value = 0.2551503493243672
$ws.Range("D25").Value = 0.08916144265243986
$ws.Range("E25").Value = 0.09297587673054908
$ws.Range("G25").Value = 0.002451087204084079
$ws.Range("I25").Value = 0.8231956217550902
$ws.Range("L25").Value = 0.2122981602508176
$ws.Range("M25").Value = 0.2708074773737152
$ws.Range("O25").Value = 3.486936415483171
